$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J27"), "https://www.adafruit.com/product/269")
$ws.Range("J27").Value2 = "https://www.adafruit.com/product/269"
